$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.004.50"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.652.80"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.14"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3903"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3827"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.37"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.351"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.001"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08450"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.86"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.078"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.009"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.647.50"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.66"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07004"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.70"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.982"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.81"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.003.90"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.441"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.984"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.11"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.39"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.424"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "138.14"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.961"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.517"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.826.74"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08091"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.737"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02935"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2679"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.72"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09128"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7608"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.41"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.03%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.32"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6968"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.466"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.104"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.000"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08332"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.85"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.223"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.96%  "
